$d = $word.ActiveDocument

# 1. Page orientation -> portrait (adds w:orient="portrait" to pgSz)
$d.PageSetup.Orientation = 0

# 2. Locate the target paragraph (last "real" paragraph of the body,
#    the one ending in "... e da frequência." with the _GoBack bookmark).
$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Contudo, durante a aula*") {
        $target = $p
        break
    }
}

# First-line indent of 708 twips (35.4 pt) in addition to existing left indent.
$target.Range.ParagraphFormat.FirstLineIndent = 35.4

# 3. Underline the word "ver" in "...gerador de função podemos ver como..."
$rng1 = $target.Range
$found1 = $rng1.Find.Execute("ver", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Font.Underline = 1
}

# 4. Underline the word "amplitude" in "...a alteração da amplitude e da..."
$rng2 = $target.Range
$found2 = $rng2.Find.Execute("amplitude", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Font.Underline = 1
}
